$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H17").Value = 1879.1072
$ws.Range("I17").Value = 1363.6364
$ws.Range("J17").Value = 2212.647
$ws.Range("K17").Value = 4090.9092
$ws.Range("L17").Value = 6637.941
$ws.Range("M17").Value = -3922.9092
$ws.Range("N17").Value = -6973.941
$ws.Range("H113").Value = 1677.25
$ws.Range("I113").Value = 1677.25
$ws.Range("K113").Value = 1677.25
$ws.Range("M113").Value = 1576.75
$ws.Range("H117").Value = 70577.336
$ws.Range("J117").Value = 70577.336
$ws.Range("L117").Value = 70577.336
$ws.Range("N117").Value = -79755.336

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2923
$ws.Range("I45").Value = 2166.5557
$ws.Range("J45").Value = 4625
$ws.Range("K45").Value = 2166.5557
$ws.Range("L45").Value = 4625
$ws.Range("M45").Value = -1789.5557
$ws.Range("N45").Value = -5379
$ws.Range("H63").Value = 3918.2222
$ws.Range("I63").Value = 1980.5
$ws.Range("J63").Value = 19420
$ws.Range("K63").Value = 1980.5
$ws.Range("L63").Value = 19420
$ws.Range("M63").Value = -1294.5
$ws.Range("N63").Value = -20792
$ws.Range("H66").Value = 3918.2222
$ws.Range("I66").Value = 1980.5
$ws.Range("J66").Value = 19420
$ws.Range("K66").Value = 9902.5
$ws.Range("L66").Value = 97100
$ws.Range("M66").Value = -6470.5
$ws.Range("N66").Value = -103964
$ws.Range("H97").Value = 1202.3077
$ws.Range("I97").Value = 651
$ws.Range("K97").Value = 651
$ws.Range("M97").Value = -155
$ws.Range("H110").Value = 3047.389
$ws.Range("I110").Value = 2575.3076
$ws.Range("J110").Value = 4274.8
$ws.Range("K110").Value = 2575.3076
$ws.Range("L110").Value = 4274.8
$ws.Range("M110").Value = -530.3076000000001
$ws.Range("N110").Value = -8364.799999999999

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 8441.111000000001
$ws.Range("I86").Value = 2995
$ws.Range("J86").Value = 9997.143
$ws.Range("K86").Value = 2995
$ws.Range("L86").Value = 9997.143
$ws.Range("M86").Value = -1872
$ws.Range("N86").Value = -12243.143
$ws.Range("H89").Value = 8441.111000000001
$ws.Range("I89").Value = 2995
$ws.Range("J89").Value = 9997.143
$ws.Range("K89").Value = 14975
$ws.Range("L89").Value = 49985.715
$ws.Range("M89").Value = -9359
$ws.Range("N89").Value = -61217.715
$ws.Range("H94").Value = 778
$ws.Range("J94").Value = 100
$ws.Range("L94").Value = 100
$ws.Range("N94").Value = -1002
$ws.Range("H103").Value = 40000
$ws.Range("J103").Value = 40000
$ws.Range("L103").Value = 40000
$ws.Range("N103").Value = -42344
$ws.Range("H105").Value = 1785.4375
$ws.Range("I105").Value = 1608.3
$ws.Range("K105").Value = 1608.3
$ws.Range("M105").Value = 138.7
$ws.Range("H107").Value = 5416.3477
$ws.Range("I107").Value = 1668.3
$ws.Range("K107").Value = 1668.3
$ws.Range("M107").Value = 251.7

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3024.818
$ws.Range("I22").Value = 2503.8333
$ws.Range("J22").Value = 3650
$ws.Range("K22").Value = 2503.8333
$ws.Range("L22").Value = 3650
$ws.Range("M22").Value = -2153.8333
$ws.Range("N22").Value = -4350

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 2220.2
$ws.Range("J35").Value = 2749.75
$ws.Range("L35").Value = 8249.25
$ws.Range("N35").Value = -8825.25
$ws.Range("H39").Value = 5421.5
$ws.Range("J39").Value = 7471.143
$ws.Range("L39").Value = 22413.429
$ws.Range("N39").Value = -23001.429
$ws.Range("H80").Value = 4912.913
$ws.Range("I80").Value = 4736.684
$ws.Range("J80").Value = 5750
$ws.Range("K80").Value = 14210.052
$ws.Range("L80").Value = 17250
$ws.Range("M80").Value = -13274.052
$ws.Range("N80").Value = -19122
$ws.Range("H83").Value = 4912.913
$ws.Range("I83").Value = 4736.684
$ws.Range("J83").Value = 5750
$ws.Range("K83").Value = 42630.156
$ws.Range("L83").Value = 51750
$ws.Range("M83").Value = -37950.156
$ws.Range("N83").Value = -61110
$ws.Range("H129").Value = 1376.7
$ws.Range("I129").Value = 825.5714
$ws.Range("J129").Value = 2662.6667
$ws.Range("K129").Value = 2476.7142
$ws.Range("L129").Value = 7988.000100000001
$ws.Range("M129").Value = 2523.2858
$ws.Range("N129").Value = -17988.0001
$ws.Range("H131").Value = 2135.7368
$ws.Range("I131").Value = 1725.4546
$ws.Range("J131").Value = 2699.875
$ws.Range("K131").Value = 5176.3638
$ws.Range("L131").Value = 8099.625
$ws.Range("M131").Value = -136.3638000000001
$ws.Range("N131").Value = -18179.625

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 500
$ws.Range("J107").Value = 500
$ws.Range("L107").Value = 500
$ws.Range("N107").Value = -4340
$ws.Range("H113").Value = 9665.5
$ws.Range("I113").Value = 8832.333000000001
$ws.Range("J113").Value = 9943.223
$ws.Range("K113").Value = 8832.333000000001
$ws.Range("L113").Value = 9943.223
$ws.Range("M113").Value = -6662.333000000001
$ws.Range("N113").Value = -14283.223
$ws.Range("H122").Value = 419843.4
$ws.Range("I122").Value = 457304
$ws.Range("K122").Value = 1371912
$ws.Range("M122").Value = -1369462
$ws.Range("H133").Value = 130780
$ws.Range("J133").Value = 130780
$ws.Range("L133").Value = 130780
$ws.Range("N133").Value = -140900

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 7401.615
$ws.Range("I100").Value = 6111
$ws.Range("K100").Value = 6111
$ws.Range("M100").Value = -5570
$ws.Range("H125").Value = 99928.75
$ws.Range("J125").Value = 99928.75
$ws.Range("L125").Value = 99928.75
$ws.Range("N125").Value = -109768.75

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2667.3333
$ws.Range("I96").Value = 2667.3333
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 2667.3333
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -1294.3333
$ws.Range("N96").ClearContents()
$ws.Range("H107").Value = 585.5714
$ws.Range("I107").Value = 499.81818
$ws.Range("K107").Value = 1499.45454
$ws.Range("M107").Value = 420.54546
$ws.Range("H122").Value = 2489
$ws.Range("I122").Value = 992.6
$ws.Range("J122").Value = 4359.5
$ws.Range("K122").Value = 2977.8
$ws.Range("L122").Value = 13078.5
$ws.Range("M122").Value = -527.8000000000002
$ws.Range("N122").Value = -17978.5
$ws.Range("H132").Value = 4137.6
$ws.Range("I132").Value = 2646
$ws.Range("J132").Value = 6375
$ws.Range("K132").Value = 7938
$ws.Range("L132").Value = 19125
$ws.Range("M132").Value = -5408
$ws.Range("N132").Value = -24185
